# Attempted some resolution of the while loop freeze, no luck so far.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in a new log entry on row 8: Feature / Status / Duration columns
$ws.Range("B8").Value = "Basic Functionality"
$ws.Range("C8").Value = "In-Progress"
$ws.Range("D8").Value = "1 hour"

# Leave the selection where the author last left it
$ws.Range("I7").Select()
